$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2, 3, 4 get rotated: row2 <- old row3, row3 <- old row4, row4 <- old row2
# for columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), S (Precio $/Kg).

$oldD2 = $ws.Range("D2").Value2
$oldM2 = $ws.Range("M2").Value2
$oldN2 = $ws.Range("N2").Value2
$oldO2 = $ws.Range("O2").Value2
$oldP2 = $ws.Range("P2").Value2
$oldS2 = $ws.Range("S2").Value2

$oldD3 = $ws.Range("D3").Value2
$oldM3 = $ws.Range("M3").Value2
$oldN3 = $ws.Range("N3").Value2
$oldO3 = $ws.Range("O3").Value2
$oldP3 = $ws.Range("P3").Value2
$oldS3 = $ws.Range("S3").Value2

$oldD4 = $ws.Range("D4").Value2
$oldM4 = $ws.Range("M4").Value2
$oldN4 = $ws.Range("N4").Value2
$oldO4 = $ws.Range("O4").Value2
$oldP4 = $ws.Range("P4").Value2
$oldS4 = $ws.Range("S4").Value2

# Row 2 <- old row 3
$ws.Range("D2").Value2 = $oldD3
$ws.Range("M2").Value2 = $oldM3
$ws.Range("N2").Value2 = $oldN3
$ws.Range("O2").Value2 = $oldO3
$ws.Range("P2").Value2 = $oldP3
$ws.Range("S2").Value2 = $oldS3

# Row 3 <- old row 4
$ws.Range("D3").Value2 = $oldD4
$ws.Range("M3").Value2 = $oldM4
$ws.Range("N3").Value2 = $oldN4
$ws.Range("O3").Value2 = $oldO4
$ws.Range("P3").Value2 = $oldP4
$ws.Range("S3").Value2 = $oldS4

# Row 4 <- old row 2
$ws.Range("D4").Value2 = $oldD2
$ws.Range("M4").Value2 = $oldM2
$ws.Range("N4").Value2 = $oldN2
$ws.Range("O4").Value2 = $oldO2
$ws.Range("P4").Value2 = $oldP2
$ws.Range("S4").Value2 = $oldS2
